$d = $word.ActiveDocument

# Locate the date-day run ("16") in the cover letter header table (e.g. "16/08/2020").
$rng = $d.Content
$found = $rng.Find.Execute("16", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the '16' date text to update"
}

$dateStart = $rng.Start
$dateEnd = $rng.End

# Move the "_GoBack" bookmark (Word re-creates this automatically at the location of the
# most recent edit) to sit right after the edited run, before the following "/" run.
# Adding a bookmark with an already-existing name relocates it, removing it from its old spot.
$bmRange = $d.Range($dateEnd, $dateEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Now update the actual text of the run from "16" to "20".
$dateRange = $d.Range($dateStart, $dateEnd)
$dateRange.Text = "20"
